# Edit script for LOB1036.xlsx per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column B width (was incorrectly narrow due to an overlapping <col> definition) ---
$ws.Columns.Item(2).ColumnWidth = 60.7109375

# --- Objetivos/Objectives text realignment (rows 10-11 keep position, text corrected) ---
$ws.Range("B10").Value = 'Fornecer fundamentos teóricos sobre vetores, retas no espaço e plano (com suas relações), cônicas e quádricas, tópicos essenciais no estudo de todas Engenharias'
$ws.Range("C10").Value = 'Fornecer fundamentos teóricos sobre vetores, retas no espaço e plano (com suas relações), cônicas e quádricas, tópicos essenciais no estudo de todas Engenharias'

# --- Insert a new row at 13 for "Docentes responsaveis" name row; shifts old rows 13-21 down to 14-22 ---
$ws.Rows.Item(13).Insert() | Out-Null

# The inserted row 13 wrongly inherits bold formatting from row 12 (label column style).
# Reset A13 to the default/no style (it stays empty) and fix B13 to the normal wrapped style.
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true

# --- Set final content + row heights for rows 13-22 ---

# Row 13
$ws.Range("B13").Value = '3682251 - Gabrielle Weber Martins'
$ws.Range("C13").Value = '3682251 - Gabrielle Weber Martins'

# Row 14
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Vetores. Vetores no R2  e no R3. Dependência Linear. Produtos de Vetores. A Reta. O Plano. Distâncias. Coordenadas Polares. Mudança de Coordenadas. Cônicas. Superfícies Quádricas. Equações Paramétricas.'
$ws.Range("C14").Value = 'Vetores. Vetores no R2  e no R3. Dependência Linear. Produtos de Vetores. A Reta. O Plano. Distâncias. Coordenadas Polares. Mudança de Coordenadas. Cônicas. Superfícies Quádricas. Equações Paramétricas.'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Vectors. Vectors in 2 and 3 Dimensions. Linear Dependence. Products of Vectors. Lines. Planes. Distances. Polar Coordinates. Coordinates changing. Conic Sections. Quadric Surfaces.'
$ws.Range("C15").Value = 'Vectors. Vectors in 2 and 3 Dimensions. Linear Dependence. Products of Vectors. Lines. Planes. Distances. Polar Coordinates. Coordinates changing. Conic Sections. Quadric Surfaces.'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '•Vetores: Reta orientada. Eixo. Segmento orientado. Segmentos equipolentes.  Vetor. Operações com vetores. Ângulo de dois vetores.•Vetores no r2 e no r3: Decomposição de um vetor no plano. Expressão analítica de um vetor. Igualdade e operações; Vetor definido pelas coordenadas da origem e da extremidade. Decomposição de um vetor no espaço. Igualdade – Operações – Vetor definido pelos pontos extremos. Condição de paralelismo de dois vetores.•Dependência linear: Dependência e Independência Linear de vetores no R2 e no R3. Base. Mudança de Base.•Produtos de vetores: Produto escalar. Módulo de um vetor. Propriedades do produto escalar. Ângulo de dois vetores. Ângulos diretores e cosenos diretores de um vetor. Projeção de um vetor. Produto escalar no R2. Produto vetorial. Propriedades do produto vetorial. Interpretação geométrica do módulo do produto vetorial de dois vetores. Produto misto. Propriedades do produto misto. Interpretação geométrica do módulo do produto misto.•A reta: Equação vetorial da reta. Reta definida por dois pontos. Equações paramétricas da reta. Equações simétricas da reta. Equações reduzidas da reta. Retas paralelas aos planos e aos eixos coordenados. Ângulo de duas retas. Condição de paralelismo e de ortogonalidade de duas retas. Condição de coplanaridade de duas retas. Posições relativas de duas retas. Reta ortogonal a duas retas. Ponto que divide um segmento de reta em uma razão dada.•O plano: Equação geral do plano. Determinação de um plano. Planos paralelos aos eixos e aos planos coordenados – Casos particulares. Equações paramétricas do plano. Ângulo de dois planos. Ângulo de uma reta com um plano. Intersecção de dois planos. Intersecção de reta com plano.•Distâncias: Distância entre dois pontos. Distância de um ponto a uma reta. Distância de duas retas. Distância de um ponto a um plano. Distância entre dois planos. Distância de uma reta a um plano.•Coordenadas polares: Definição de Coordenadas polares, equações e gráficos polares. Relacionando coordenadas polares e coordenadas cartesianas•Mudança de coordenadas: Mudança de coordenadas em R2 e em R3. Aplicação de translações e rotações.•Equações paramétricas: da reta, da circunferência. Equações Paramétricas de curvas.•Cônicas: A parábola. A elipse. A hipérbole. As seções cônicas.•Superfícies quádricas: Introdução. Superfícies quádricas centradas. Superfícies quádricas não'
$ws.Range("C16").Value = '•Vetores: Reta orientada. Eixo. Segmento orientado. Segmentos equipolentes.  Vetor. Operações com vetores. Ângulo de dois vetores.•Vetores no r2 e no r3: Decomposição de um vetor no plano. Expressão analítica de um vetor. Igualdade e operações; Vetor definido pelas coordenadas da origem e da extremidade. Decomposição de um vetor no espaço. Igualdade – Operações – Vetor definido pelos pontos extremos. Condição de paralelismo de dois vetores.•Dependência linear: Dependência e Independência Linear de vetores no R2 e no R3. Base. Mudança de Base.•Produtos de vetores: Produto escalar. Módulo de um vetor. Propriedades do produto escalar. Ângulo de dois vetores. Ângulos diretores e cosenos diretores de um vetor. Projeção de um vetor. Produto escalar no R2. Produto vetorial. Propriedades do produto vetorial. Interpretação geométrica do módulo do produto vetorial de dois vetores. Produto misto. Propriedades do produto misto. Interpretação geométrica do módulo do produto misto.•A reta: Equação vetorial da reta. Reta definida por dois pontos. Equações paramétricas da reta. Equações simétricas da reta. Equações reduzidas da reta. Retas paralelas aos planos e aos eixos coordenados. Ângulo de duas retas. Condição de paralelismo e de ortogonalidade de duas retas. Condição de coplanaridade de duas retas. Posições relativas de duas retas. Reta ortogonal a duas retas. Ponto que divide um segmento de reta em uma razão dada.•O plano: Equação geral do plano. Determinação de um plano. Planos paralelos aos eixos e aos planos coordenados – Casos particulares. Equações paramétricas do plano. Ângulo de dois planos. Ângulo de uma reta com um plano. Intersecção de dois planos. Intersecção de reta com plano.•Distâncias: Distância entre dois pontos. Distância de um ponto a uma reta. Distância de duas retas. Distância de um ponto a um plano. Distância entre dois planos. Distância de uma reta a um plano.•Coordenadas polares: Definição de Coordenadas polares, equações e gráficos polares. Relacionando coordenadas polares e coordenadas cartesianas•Mudança de coordenadas: Mudança de coordenadas em R2 e em R3. Aplicação de translações e rotações.•Equações paramétricas: da reta, da circunferência. Equações Paramétricas de curvas.•Cônicas: A parábola. A elipse. A hipérbole. As seções cônicas.•Superfícies quádricas: Introdução. Superfícies quádricas centradas. Superfícies quádricas não'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '•Vectors. Directed line. Axis. Directed line segment. Equipollent line segments. Vector. Addition and scalar multiples ofvectors. Angle between two vectors.•Vectors in 2 and 3 dimensions: Decomposition of a vector in 2 dimensions. Analytical expression of a vector. Equal vectors and Addition and scalar multiples of vectors. Vector defined by initial and terminal points coordinates. Decomposition of a vector in 3 dimensions. Equal vectors and Addition and scalar multiples of vectors. Vector defined by extreme points. Parallel vectors.•Linear dependence: Vectors linear dependence and interdependence in 2 and 3 dimensions. Base. Base changing.•Products of vectors. Dot product. Magnitude of a vector. Properties of the dot product. Angle between two vectors. Angles directors and cosines directors of a vector. Vector projection. Dot product in 2 dimensions. Cross product. Properties of the cross product. Geometric interpretation of the cross product magnitude. Scalar triple roduct. Properties of the scalar tripleproduct. Geometric interpretation of the scalar triple roduct magnitude.•Lines. Vector equation of the line. Line defined by two points. Parametric equations of the line. Symmetric equations of the line. Reduced equations of the line. Lines parallel to the coordinate planes and axes. Angle between two lines. Parallel and orthogonal lines. Coplanar lines. Relative positions between two lines. Line orthogonal to two lines. Point that divides a line segment in a given ratio.•Planes. Linear equation of the plane. Determination of a plane. Planes parallel to the coordinate planes and axes. Parametric equation of the plane. Angle between two planes. Angle between a line and a plane. Intersection of two planes. Intersection of a line and a plane.•Distances. Distance between two points. Distance from a point to a line. Distance between two lines. Distance from a point to a plane. Distance between two planes. Distance from a line to a plane.•Polar coordinates. Polar coordinates definition, polar equations and graphics. Relating polar coordinates to Cartesian coordinates.•Coordinate changing: Changing coordinates in 2 and 3 dimensions. Translatory and rotations applications.•Parametric equations: line, circumference. Parametric equations for curves.•Conic sections. The parabola. The ellipse. The hyperbola. The conic sections.•Quadric surfaces. Introduction. Centered quadric surfaces. Noncentered quadric surfaces. Cones. Cylinders.'
$ws.Range("C17").Value = '•Vectors. Directed line. Axis. Directed line segment. Equipollent line segments. Vector. Addition and scalar multiples ofvectors. Angle between two vectors.•Vectors in 2 and 3 dimensions: Decomposition of a vector in 2 dimensions. Analytical expression of a vector. Equal vectors and Addition and scalar multiples of vectors. Vector defined by initial and terminal points coordinates. Decomposition of a vector in 3 dimensions. Equal vectors and Addition and scalar multiples of vectors. Vector defined by extreme points. Parallel vectors.•Linear dependence: Vectors linear dependence and interdependence in 2 and 3 dimensions. Base. Base changing.•Products of vectors. Dot product. Magnitude of a vector. Properties of the dot product. Angle between two vectors. Angles directors and cosines directors of a vector. Vector projection. Dot product in 2 dimensions. Cross product. Properties of the cross product. Geometric interpretation of the cross product magnitude. Scalar triple roduct. Properties of the scalar tripleproduct. Geometric interpretation of the scalar triple roduct magnitude.•Lines. Vector equation of the line. Line defined by two points. Parametric equations of the line. Symmetric equations of the line. Reduced equations of the line. Lines parallel to the coordinate planes and axes. Angle between two lines. Parallel and orthogonal lines. Coplanar lines. Relative positions between two lines. Line orthogonal to two lines. Point that divides a line segment in a given ratio.•Planes. Linear equation of the plane. Determination of a plane. Planes parallel to the coordinate planes and axes. Parametric equation of the plane. Angle between two planes. Angle between a line and a plane. Intersection of two planes. Intersection of a line and a plane.•Distances. Distance between two points. Distance from a point to a line. Distance between two lines. Distance from a point to a plane. Distance between two planes. Distance from a line to a plane.•Polar coordinates. Polar coordinates definition, polar equations and graphics. Relating polar coordinates to Cartesian coordinates.•Coordinate changing: Changing coordinates in 2 and 3 dimensions. Translatory and rotations applications.•Parametric equations: line, circumference. Parametric equations for curves.•Conic sections. The parabola. The ellipse. The hyperbola. The conic sections.•Quadric surfaces. Introduction. Centered quadric surfaces. Noncentered quadric surfaces. Cones. Cylinders.'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A18").Value = 'Avaliação:'

# Row 19
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = '1.CAMARGO, Ivan ; BOULOS, Paulo. Geometria Analítica: um tratamento vetorial. São Paulo: Prentice Hall, 2005.2.LIMA, Elon Lages de. Geometria analítica e algebra Linear. Rio de Janeiro: SBM SociedadeBrasileira de Matemática,2001. Coleção Matemática Universitária.3.CAROLI, Alésio de; CALLIOLI, A.; FEITOSA, Miguel O. Matrizes vetores geometria analítica. São Paulo: Nobel, 1998.4.SANTOS, Nathan Moreira dos. Vetores e matrizes: uma introdução à álgebra linear. São Paulo: Thomson, 2007.'
$ws.Range("C22").Value = '1.CAMARGO, Ivan ; BOULOS, Paulo. Geometria Analítica: um tratamento vetorial. São Paulo: Prentice Hall, 2005.2.LIMA, Elon Lages de. Geometria analítica e algebra Linear. Rio de Janeiro: SBM SociedadeBrasileira de Matemática,2001. Coleção Matemática Universitária.3.CAROLI, Alésio de; CALLIOLI, A.; FEITOSA, Miguel O. Matrizes vetores geometria analítica. São Paulo: Nobel, 1998.4.SANTOS, Nathan Moreira dos. Vetores e matrizes: uma introdução à álgebra linear. São Paulo: Thomson, 2007.'
$ws.Rows.Item(22).RowHeight = 120
